# 🚌 141: 30/12 14:21 LP1912+6203+6173
# Appends freshly-scraped rows to the "LP1912" and "6203-6173" sheets and
# refreshes the "Última actualización" timestamp (and "Total filas" counts)
# on all three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912  (columns: A info, B Hora_Scrap, C Hora_Llegada, D Linea,
#                    E Minutos, F Parada, G Fecha)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 11:21:04"
$ws1.Cells.Item(3,1).Value = "Total filas: 170"

$sheet1Rows = @(
  @(154,"11:20:54","11:25","10_OLMOS",5,"LP1912","30/12/2025"),
  @(155,"11:20:54","11:26","16_SANTA ANA",6,"LP1912","30/12/2025"),
  @(156,"11:20:54","11:26","225_C ROCA-H SUR",6,"LP1912","30/12/2025"),
  @(157,"11:20:54","11:34","23_HERNANDEZ",14,"LP1912","30/12/2025"),
  @(158,"11:20:54","11:41","17_ROMERO",21,"LP1912","30/12/2025"),
  @(159,"11:20:54","11:43","10_OLMOS",23,"LP1912","30/12/2025"),
  @(160,"11:20:54","11:52","15_ABASTO",32,"LP1912","30/12/2025"),
  @(161,"11:20:54","12:02","84_COLONIA URQUIZA-ESC 49",42,"LP1912","30/12/2025"),
  @(162,"11:20:54","12:04","23_HERNANDEZ",44,"LP1912","30/12/2025"),
  @(163,"11:20:54","12:06","16_P MOR-SANTA ANA",46,"LP1912","30/12/2025"),
  @(164,"11:20:54","12:20","14_ABASTO",60,"LP1912","30/12/2025"),
  @(165,"11:20:54","12:20","26_HERNANDEZ",60,"LP1912","30/12/2025"),
  @(166,"11:20:54","12:34","23_HERNANDEZ",74,"LP1912","30/12/2025"),
  @(167,"11:20:54","12:38","17_179 Y 38",78,"LP1912","30/12/2025"),
  @(168,"11:20:54","12:48","11_ETCHEVERRY",88,"LP1912","30/12/2025"),
  @(169,"11:20:54","12:50","17_ROMERO",90,"LP1912","30/12/2025"),
  @(170,"11:20:54","12:54","10_OLMOS",94,"LP1912","30/12/2025"),
  @(171,"11:20:54","12:55","15_ABASTO",95,"LP1912","30/12/2025")
)

foreach ($row in $sheet1Rows) {
  $r = $row[0]
  $ws1.Cells.Item($r,1).Value = ""
  $ws1.Cells.Item($r,2).Value = $row[1]
  $ws1.Cells.Item($r,3).Value = $row[2]
  $ws1.Cells.Item($r,4).Value = $row[3]
  $ws1.Cells.Item($r,5).Value = $row[4]
  $ws1.Cells.Item($r,6).Value = $row[5]
  $ws1.Cells.Item($r,7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215  -- only the timestamp refreshes, no new rows.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 11:21:04"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173  (columns: A info, B Fecha, C Hora_Scrap, D Hora_Llegada,
#                       E Linea, F Minutos, G Parada)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 11:21:04"
$ws3.Cells.Item(3,1).Value = "Total filas: 23"

$sheet3Rows = @(
  @(23,"30/12/2025","11:21:04","12:04","215A_LA PLATA",43,"L6173"),
  @(24,"30/12/2025","11:20:59","12:53","215C_LA PLATA",93,"L6203")
)

foreach ($row in $sheet3Rows) {
  $r = $row[0]
  $ws3.Cells.Item($r,1).Value = ""
  $ws3.Cells.Item($r,2).Value = $row[1]
  $ws3.Cells.Item($r,3).Value = $row[2]
  $ws3.Cells.Item($r,4).Value = $row[3]
  $ws3.Cells.Item($r,5).Value = $row[4]
  $ws3.Cells.Item($r,6).Value = $row[5]
  $ws3.Cells.Item($r,7).Value = $row[6]
}
